# "new layout with three 2d views"
# Adds a second monthly time-tracking block (rows 16-28), mirroring the
# first block (rows 2-15), and renames the "@Parsiss" label to "@IACT".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant used throughout to copy a source cell's number
# format / font / fill / border / alignment onto a just-written cell.
$xlPasteFormats = -4122

function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------
# 1) Rename the "@Parsiss" entry (row 14) to "@IACT"
# ---------------------------------------------------------------------
$ws.Range("C14").Value = "@IACT"
Copy-Format "C15" "C14"

# ---------------------------------------------------------------------
# 2) New separator rows below the first block
# ---------------------------------------------------------------------
# Row 16 mirrors the thick-bottom separator of row 3 (only F16:G16 here)
Copy-Format "B3" "F16"
Copy-Format "B3" "G16"
$ws.Rows(16).RowHeight = 15

# Row 17 mirrors the blank spacer row 4
foreach ($col in @("A", "B", "C", "D", "E")) {
    Copy-Format ($col + "4") ($col + "17")
}

# ---------------------------------------------------------------------
# 3) Second table header (row 18), mirrors row 5
# ---------------------------------------------------------------------
$ws.Range("A18").Value = "دی 99"
Copy-Format "A5" "A18"
$ws.Range("B18").Value = "Activity"
Copy-Format "B5" "B18"
$ws.Range("C18").Value = "Hours"
Copy-Format "C5" "C18"
$ws.Range("E18").Value = "Tasks Done"
Copy-Format "E5" "E18"

# ---------------------------------------------------------------------
# 4) Second table body (rows 19-25), mirrors rows 6-12
# ---------------------------------------------------------------------
$ws.Range("B19").Value = "* Segmentation"
Copy-Format "B6" "B19"
$ws.Range("C19").Value = 2
Copy-Format "C6" "C19"
$ws.Range("E19").Value = "• Presentations @IACT"
Copy-Format "E6" "E19"

$ws.Range("B20").Value = "* Debug & Refactor"
Copy-Format "B7" "B20"
$ws.Range("C20").Value = 1
Copy-Format "C7" "C20"
$ws.Range("E20").Value = "• More realistic virtual view"
Copy-Format "E7" "E20"

$ws.Range("B21").Value = "* Registration"
Copy-Format "B8" "B21"
$ws.Range("C21").Value = 2
Copy-Format "C8" "C21"
$ws.Range("E21").Value = "• BronchoVision Video"
Copy-Format "E8" "E21"

$ws.Range("B22").Value = "* Tracker"
Copy-Format "B9" "B22"
$ws.Range("C22").Value = 5
Copy-Format "C9" "C22"
$ws.Range("E22").Value = "• Bug and Crash Fixes"
Copy-Format "E9" "E22"

$ws.Range("B23").Value = "* 2D/3D Views"
Copy-Format "B10" "B23"
$ws.Range("C23").Value = 3
Copy-Format "C10" "C23"

$ws.Range("B24").Value = "* BronchoVision Video"
Copy-Format "B12" "B24"
$ws.Range("C24").Value = 4
Copy-Format "C12" "C24"
Copy-Format "E12" "E24"
$ws.Rows(24).RowHeight = 15.6

$ws.Range("B25").Value = "* Meetings & Presentations"
Copy-Format "B11" "B25"
$ws.Range("C25").Value = 28
Copy-Format "C11" "C25"

# ---------------------------------------------------------------------
# 5) Totals section (rows 26-28), mirrors rows 13-15
# ---------------------------------------------------------------------
$ws.Range("B26").Value = "• Total Hours"
Copy-Format "B13" "B26"
$ws.Range("C26").Formula = "=SUM(C19:C25)"
Copy-Format "C13" "C26"

$ws.Range("C27").Value = "@IACT"
Copy-Format "C14" "C27"
$ws.Range("D27").Value = 32
Copy-Format "D14" "D27"

$ws.Range("C28").Value = "@Home"
Copy-Format "C15" "C28"
$ws.Range("D28").Formula = "=C26-D27"
Copy-Format "D15" "D28"

# ---------------------------------------------------------------------
# 6) Restore the active selection to match the new layout
# ---------------------------------------------------------------------
$ws.Range("G24").Select()
